# Updated cryptos list on Tue Jun 11 19:47:52 UTC 2024 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on Sheet1 with the latest
# scraped values. Price values that look like plain decimal numbers are
# forced to text (matching the source data, which stores them as strings)
# by temporarily switching the cell to a text number format and then
# clearing that temporary formatting so the cell keeps its original style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.305.21"
$ws.Range("E2").Value = "  -3.45%  "
$ws.Range("D3").Value = "3.491.63"
$ws.Range("E3").Value = "  -5.05%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.81%  "
$ws.Range("D7").Value = "3.491.67"
$ws.Range("E7").Value = "  -4.99%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000216"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "4.080.42"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.40"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "3.488.24"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("D17").Value = "67.215.04"
$ws.Range("E17").Value = "  -3.60%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.00"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "446.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -12.75%  "
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.09"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000127"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "3.629.30"
$ws.Range("E27").Value = "  -5.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.74"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.66%  "
$ws.Range("E35").Value = "  -4.93%  "
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").Value = "3.483.89"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.98"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.65%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0873"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  -7.50%  "
$ws.Range("E45").Value = "  -5.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.81%  "
$ws.Range("E49").Value = "  -8.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("E51").Value = "  -3.73%  "
